# Apply the LinuxForHealth re-brand / refresh edit to the
# StructureDefinition-family-income-level workbook.
#
# Changes:
#  - Metadata sheet: URL, Version, Date and Publisher values updated.
#  - Elements sheet: the "Extension" summary row (row 2) no longer repeats
#    the ele-1/ext-1 constraint text in its Constraint(s) column (AI2) -
#    that text remains only on the Extension.extension row (row 4, AI4).
#  - Elements sheet: the Fixed Value for Extension.url (Q5) is refreshed to
#    mirror the new canonical URL used on the Metadata sheet.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/family-income-level"
$metadata.Range("B3").Value = "8.0.0"
$metadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$metadata.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/family-income-level"
